$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "Pricetobuy" column (G), which will
# hold a computed "Full Name" column. Excel shifts formulas/styles for us.
$ws.Columns("G:G").Insert()

# Match the width used in the target workbook for the new "Full Name" column.
$ws.Columns("G:G").ColumnWidth = 97.736979166667

# Header for the new column.
$ws.Range("G1").Value = "Full Name"

# Fill in the concatenated "Full Name" formula for every data row (2-18).
# G2 is written on its own (matches the pre-existing "Cost" column pattern,
# whose formula in row 2 also stands alone), then G3:G18 is filled as one
# block so Excel groups it into a single shared formula - mirroring how the
# original workbook already grouped the "Cost" formulas below row 2.
$ws.Range("G2").Formula = '=CONCATENATE(F2," - ",E2," - ",D2," - ",C2," - ",B2," - ",A2)'
$ws.Range("G3:G18").Formula = '=CONCATENATE(F3," - ",E3," - ",D3," - ",C3," - ",B3," - ",A3)'

# Re-assert the "Cost" column (now J) formula over its original shared-
# formula groupings (J3:J6, J7:J8, J9:J18) so the grouping metadata survives
# the column insert exactly as it was before (insert alone flattens it).
$ws.Range("J3:J6").Formula = '=H3/(I3/1.081)'
$ws.Range("J7:J8").Formula = '=H7/(I7/1.081)'
$ws.Range("J9:J18").Formula = '=H9/(I9/1.081)'

# "Updated" column (now L after the insert) gets a date stamp on every row.
# Rows 2-8 already carry the date number format (style 5) from the shifted
# column; rows 9-18 never had an "Updated" cell before, so pull the date
# format onto them too (PasteSpecial reuses the existing style instead of
# minting a new one) before writing the serial date value (2025-02-06).
$ws.Range("L2").Copy()
$ws.Range("L9:L18").PasteSpecial(-4122)
$ws.Range("L2:L18").Value = 45694

# Keep the previously selected cell sane relative to the new layout.
$ws.Range("E11").Select()
